$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("G2").Value = 6264
$ws.Range("K2").Value = 2731
$ws.Range("K3").Value = 2641
$ws.Range("E4").Value = 2031
$ws.Range("J4").Value = 1817
$ws.Range("K4").Value = 551
$ws.Range("K5").Value = 175
$ws.Range("K6").Value = 3287
$ws.Range("E7").Value = 26036
$ws.Range("G7").Value = 24709
$ws.Range("J7").Value = 29286
$ws.Range("K7").Value = 9385

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("K2").Value = 69
$ws.Range("K6").Value = 78
$ws.Range("K7").Value = 278
$ws.Range("K8").Value = 618
$ws.Range("K11").Value = 199
$ws.Range("K12").Value = 15
$ws.Range("K15").Value = 94
$ws.Range("K16").Value = 31
$ws.Range("K18").Value = 63
$ws.Range("K19").Value = 276
$ws.Range("K20").Value = 216
$ws.Range("K22").Value = 29
$ws.Range("K24").Value = 34
$ws.Range("K25").Value = 37
$ws.Range("K29").Value = 489
$ws.Range("K33").Value = 368
$ws.Range("K36").Value = 108
$ws.Range("K37").Value = 308
$ws.Range("K38").Value = 10
$ws.Range("K42").Value = 331
$ws.Range("K48").Value = 113
$ws.Range("K51").Value = 105
$ws.Range("K52").Value = 261
$ws.Range("K53").Value = 137
$ws.Range("J54").Value = 573
$ws.Range("E63").Value = 367
$ws.Range("G63").Value = 284
$ws.Range("J63").Value = 99
$ws.Range("K63").Value = 34
$ws.Range("K65").Value = 222
$ws.Range("K67").Value = 362
$ws.Range("K71").Value = 30
$ws.Range("K73").Value = 92
$ws.Range("K76").Value = 139
$ws.Range("K77").Value = 66
$ws.Range("K78").Value = 131
$ws.Range("K79").Value = 238
$ws.Range("K80").Value = 31
$ws.Range("K83").Value = 204
$ws.Range("K85").Value = 448
$ws.Range("J86").Value = 177
$ws.Range("K86").Value = 58
$ws.Range("K88").Value = 106
$ws.Range("K89").Value = 123
$ws.Range("K94").Value = 112
$ws.Range("K95").Value = 153
$ws.Range("K96").Value = 130
$ws.Range("K99").Value = 168
$ws.Range("E101").Value = 26036
$ws.Range("G101").Value = 24709
$ws.Range("J101").Value = 29286
$ws.Range("K101").Value = 9385

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("K3").Value = 19
$ws.Range("K7").Value = 130

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("K3").Value = 87
$ws.Range("K7").Value = 278

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("K2").Value = 61
$ws.Range("K3").Value = 52
$ws.Range("K6").Value = 77
$ws.Range("K7").Value = 199

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("K3").Value = 38
$ws.Range("K4").Value = 19
$ws.Range("K7").Value = 123

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("K3").Value = 154
$ws.Range("K6").Value = 100
$ws.Range("K7").Value = 448

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("K3").Value = 62
$ws.Range("K6").Value = 109
$ws.Range("K7").Value = 261

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("K2").Value = 29
$ws.Range("K7").Value = 137

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("K3").Value = 183
$ws.Range("K6").Value = 209
$ws.Range("K7").Value = 618

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("K2").Value = 77
$ws.Range("K7").Value = 204

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("K3").Value = 135
$ws.Range("K6").Value = 104
$ws.Range("K7").Value = 368

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("K2").Value = 47
$ws.Range("K3").Value = 53
$ws.Range("K7").Value = 153

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("K3").Value = 104
$ws.Range("K6").Value = 97
$ws.Range("K7").Value = 308

$ws = $wb.Worksheets.Item('New City')
$ws.Range("K3").Value = 53
$ws.Range("K6").Value = 92
$ws.Range("K7").Value = 222

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("K3").Value = 64
$ws.Range("K7").Value = 168

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("K3").Value = 114
$ws.Range("K6").Value = 106
$ws.Range("K7").Value = 362

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("J4").Value = 48
$ws.Range("J7").Value = 573

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("K2").Value = 129
$ws.Range("K3").Value = 163
$ws.Range("K5").Value = 12
$ws.Range("K6").Value = 157
$ws.Range("K7").Value = 489

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("K3").Value = 21
$ws.Range("K7").Value = 113

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("K2").Value = 88
$ws.Range("K3").Value = 72
$ws.Range("K4").Value = 11
$ws.Range("K7").Value = 276

$ws = $wb.Worksheets.Item('River North')
$ws.Range("K6").Value = 84
$ws.Range("K7").Value = 139

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("K2").Value = 27
$ws.Range("K3").Value = 27
$ws.Range("K7").Value = 78

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("K2").Value = 84
$ws.Range("K3").Value = 105
$ws.Range("K4").Value = 13
$ws.Range("K6").Value = 127
$ws.Range("K7").Value = 331

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("K3").Value = 30
$ws.Range("K6").Value = 50
$ws.Range("K7").Value = 131

$ws = $wb.Worksheets.Item('Dunning')
$ws.Range("K2").Value = 11
$ws.Range("K7").Value = 34

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("K2").Value = 78
$ws.Range("K3").Value = 86
$ws.Range("K7").Value = 238

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("K2").Value = 76
$ws.Range("K6").Value = 76
$ws.Range("K7").Value = 216

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("K6").Value = 15
$ws.Range("K7").Value = 63

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("K4").Value = 9
$ws.Range("K7").Value = 108

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("K2").Value = 31
$ws.Range("K4").Value = 11
$ws.Range("K7").Value = 112

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("K6").Value = 4
$ws.Range("K7").Value = 37

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("K4").Value = 4
$ws.Range("K6").Value = 15

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("K2").Value = 33
$ws.Range("K7").Value = 94

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("K2").Value = 28
$ws.Range("K7").Value = 92

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("K2").Value = 20
$ws.Range("K7").Value = 69

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("K6").Value = 55
$ws.Range("K7").Value = 106

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("K2").Value = 15
$ws.Range("J4").Value = 96
$ws.Range("J7").Value = 177
$ws.Range("K7").Value = 58

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("K6").Value = 39
$ws.Range("K7").Value = 105

$ws = $wb.Worksheets.Item('Clearing')
$ws.Range("K3").Value = 9
$ws.Range("K7").Value = 29

$ws = $wb.Worksheets.Item('Oakland')
$ws.Range("K6").Value = 9
$ws.Range("K7").Value = 30

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("K2").Value = 33
$ws.Range("K7").Value = 66

$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range("K2").Value = 8
$ws.Range("K7").Value = 31

$ws = $wb.Worksheets.Item('Beverly')
$ws.Range("K2").Value = 7
$ws.Range("K7").Value = 15

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range("K6").Value = 20
$ws.Range("K7").Value = 31

$ws = $wb.Worksheets.Item('Grant Park')
$ws.Range("K3").Value = 4
$ws.Range("K6").Value = 10
